$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells are forced to Text format ("@") before assignment so that
# numeric-looking strings (prices, percentages) are preserved verbatim as text,
# matching the original inlineStr string storage in the workbook.
$updates = @(
    @{ Cell = 'D2'; Value = '303.11' }
    @{ Cell = 'E2'; Value = '-3.78%' }
    @{ Cell = 'D3'; Value = '35.62' }
    @{ Cell = 'E3'; Value = '1.07%' }
    @{ Cell = 'D4'; Value = '5.054' }
    @{ Cell = 'E4'; Value = '-1.46%' }
    @{ Cell = 'D5'; Value = '0.08041' }
    @{ Cell = 'E5'; Value = '-2.03%' }
    @{ Cell = 'D6'; Value = '1.947' }
    @{ Cell = 'E6'; Value = '-8.65%' }
    @{ Cell = 'D7'; Value = '7.804' }
    @{ Cell = 'E7'; Value = '-1.99%' }
    @{ Cell = 'B8'; Value = 'BTSEToken' }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse' }
    @{ Cell = 'D8'; Value = '2.985' }
    @{ Cell = 'E8'; Value = '0.76%' }
    @{ Cell = 'B9'; Value = 'MXToken' }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'D9'; Value = '0.9268' }
    @{ Cell = 'E9'; Value = '-0.22%' }
    @{ Cell = 'B10'; Value = 'LiechtensteinCryptoassetsExchange' }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' }
    @{ Cell = 'D10'; Value = '0.1283' }
    @{ Cell = 'E10'; Value = '23.60%' }
    @{ Cell = 'B11'; Value = 'WazirX' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' }
    @{ Cell = 'D11'; Value = '0.1860' }
    @{ Cell = 'E11'; Value = '-2.04%' }
    @{ Cell = 'B12'; Value = 'MandalaExchangeToken' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' }
    @{ Cell = 'D12'; Value = '0.09420' }
    @{ Cell = 'E12'; Value = '3.10%' }
    @{ Cell = 'B13'; Value = 'BitrueCoin' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' }
    @{ Cell = 'D13'; Value = '0.03427' }
    @{ Cell = 'E13'; Value = '-5.30%' }
    @{ Cell = 'B14'; Value = 'BitMartToken' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' }
    @{ Cell = 'D14'; Value = '0.09879' }
    @{ Cell = 'E14'; Value = '-0.43%' }
    @{ Cell = 'B15'; Value = 'BitForexToken' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' }
    @{ Cell = 'D15'; Value = '0.001396' }
    @{ Cell = 'E15'; Value = '-3.30%' }
    @{ Cell = 'B16'; Value = 'TigerCash' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = 'D16'; Value = '0.005720' }
    @{ Cell = 'E16'; Value = '-0.50%' }
    @{ Cell = 'B17'; Value = 'LEO' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = 'D17'; Value = '3.514' }
    @{ Cell = 'E17'; Value = '1.36%' }
    @{ Cell = 'B18'; Value = 'GateToken' }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = 'D18'; Value = '4.060' }
    @{ Cell = 'E18'; Value = '-2.17%' }
    @{ Cell = 'D19'; Value = '0.3403' }
    @{ Cell = 'E19'; Value = '-0.92%' }
    @{ Cell = 'D20'; Value = '0.1301' }
    @{ Cell = 'E20'; Value = '-0.71%' }
    @{ Cell = 'D21'; Value = '5.063' }
    @{ Cell = 'E21'; Value = '-0.65%' }
    @{ Cell = 'D22'; Value = '0.2467' }
    @{ Cell = 'E22'; Value = '11.56%' }
    @{ Cell = 'D23'; Value = '0.04490' }
    @{ Cell = 'E23'; Value = '-1.23%' }
    @{ Cell = 'D24'; Value = '0.001215' }
    @{ Cell = 'E24'; Value = '-2.43%' }
    @{ Cell = 'D25'; Value = '0.004817' }
    @{ Cell = 'E25'; Value = '2.25%' }
    @{ Cell = 'D26'; Value = '0.0001250' }
    @{ Cell = 'E26'; Value = '0.02%' }
    @{ Cell = 'D27'; Value = '0.0002999' }
    @{ Cell = 'E27'; Value = '-33.35%' }
    @{ Cell = 'D39'; Value = '0.01915' }
    @{ Cell = 'E39'; Value = '-2.55%' }
    @{ Cell = 'D40'; Value = '0.04751' }
    @{ Cell = 'E40'; Value = '-3.52%' }
    @{ Cell = 'D41'; Value = '0.007370' }
    @{ Cell = 'E41'; Value = '-3.71%' }
    @{ Cell = 'D42'; Value = '0.009640' }
    @{ Cell = 'E42'; Value = '22.62%' }
    @{ Cell = 'D43'; Value = '0.1336' }
    @{ Cell = 'E43'; Value = '-4.11%' }
    @{ Cell = 'D44'; Value = '0.002110' }
    @{ Cell = 'E44'; Value = '0.73%' }
    @{ Cell = 'D45'; Value = '0.01087' }
    @{ Cell = 'E45'; Value = '-7.92%' }
    @{ Cell = 'D46'; Value = '0.00006259' }
    @{ Cell = 'D47'; Value = '0.00000000750' }
    @{ Cell = 'E47'; Value = '0.01%' }
    @{ Cell = 'E48'; Value = '67.38%' }
    @{ Cell = 'E49'; Value = '-12.38%' }
    @{ Cell = 'D50'; Value = '0.00002100' }
    @{ Cell = 'E50'; Value = '0.01%' }
    @{ Cell = 'D51'; Value = '0.0002000' }
    @{ Cell = 'E51'; Value = '0.01%' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}

